$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 37, shifting existing rows 37-53 down to 38-54.
$ws.Rows("37").Insert()

# Populate the newly inserted row 37 with the new data record.
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44510
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = "Otros"
$ws.Range("I37").Value = 100107002
$ws.Range("J37").Value = "Chirimoya"
$ws.Range("K37").Value = "Cultivar IV Región"
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 150
$ws.Range("N37").Value = 25000
$ws.Range("O37").Value = 25000
$ws.Range("P37").Value = 25000
$ws.Range("Q37").Value = "`$/bandeja 10 kilos"
$ws.Range("R37").Value = "Provincia de Limarí"
$ws.Range("S37").Value = 2500
$ws.Range("T37").Value = 10
